# Update 2018-05-18#8 - add ErrorMessages section (SaveErrorMsgML) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- New "ErrorMessages" section header row (merged, bold, left-aligned,
#     matching the style used by the other section headers in the sheet) ---
$ws.Range("A234").Value = "ErrorMessages"
$ws.Range("A234").Font.Bold = $true
$ws.Range("B234").Font.Bold = $true
$ws.Range("C234").Font.Bold = $true
$ws.Range("A234:C234").HorizontalAlignment = -4131
$ws.Range("A234:C234").Merge()

# --- New data rows for the ErrorMessages / SaveErrorMsgML stored proc.
#     Column C ("SaveErrorMsgML") is written before column B on each row so
#     that the shared-string table picks up new unique strings in the same
#     order the source workbook used (ErrorMessages, SaveErrorMsgML,
#     "Error Code cannot be null...", "Error Message (ML) cannot be null..."). ---
$ws.Range("A235").Value = 2201
$ws.Range("C235").Value = "SaveErrorMsgML"
$ws.Range("B235").Value = "Error Code cannot be null or empty string."

$ws.Range("A236").Value = 2202
$ws.Range("C236").Value = "SaveErrorMsgML"
$ws.Range("B236").Value = "Language Id cannot be null or empty string."

$ws.Range("A237").Value = 2203
$ws.Range("C237").Value = "SaveErrorMsgML"
$ws.Range("B237").Value = "Language Id not found."

$ws.Range("A238").Value = 2204
$ws.Range("C238").Value = "SaveErrorMsgML"
$ws.Range("B238").Value = "Error Message (ML) cannot be null or empty string."

# --- Update the saved view state (scroll position / active selection) to
#     match where the workbook was left after the edit. ---
$excel.ActiveWindow.ScrollRow = 222
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B241").Select() | Out-Null
